# Add two new data sheets ("census" and "proportion_calf") holding data that
# used to be passed as bpt_analyse() function arguments.
#
# Final tab order: location, event, census, proportion_calf
# "census" becomes the active/selected sheet (as in the authored workbook).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Create the two new worksheets, after "event", in the order they must end
# up in the tab strip.
$wsCensus = $wb.Worksheets.Add($null, $lastSheet)
$wsCensus.Name = "census"

$wsCalf = $wb.Worksheets.Add($null, $wsCensus)
$wsCalf.Name = "proportion_calf"

# ---------------------------------------------------------------------
# Populate "proportion_calf" FIRST so its header strings land in the
# shared-string table before "census" - matching the source workbook's
# shared string ordering (proportion_calf_* before census_*).
# ---------------------------------------------------------------------
$wsCalf.Range("A1").Value = "proportion_calf_year"
$wsCalf.Range("B1").Value = "proportion_calf_month"
$wsCalf.Range("C1").Value = "proportion_calf_day"
$wsCalf.Range("D1").Value = "proportion_calf"
$wsCalf.Range("E1").Value = "proportion_calf_cv"

$wsCalf.Range("A2").Value = 2021
$wsCalf.Range("B2").Value = 3
$wsCalf.Range("C2").Value = 31
$wsCalf.Range("D2").Value = 0.2
$wsCalf.Range("E2").Value = 0.05

$wsCalf.Range("A3").Value = 2022
$wsCalf.Range("B3").Value = 3
$wsCalf.Range("C3").Value = 31
$wsCalf.Range("D3").Value = 0.15
$wsCalf.Range("E3").Value = 0.09

# ---------------------------------------------------------------------
# Populate "census"
# ---------------------------------------------------------------------
$wsCensus.Range("A1").Value = "census_year"
$wsCensus.Range("B1").Value = "census_month"
$wsCensus.Range("C1").Value = "census_day"
$wsCensus.Range("D1").Value = "census"
$wsCensus.Range("E1").Value = "census_cv"

$wsCensus.Range("A2").Value = 2021
$wsCensus.Range("B2").Value = 3
$wsCensus.Range("C2").Value = 31
$wsCensus.Range("D2").Value = 250
$wsCensus.Range("E2").Value = 0.05

$wsCensus.Range("A3").Value = 2022
$wsCensus.Range("B3").Value = 3
$wsCensus.Range("C3").Value = 31
$wsCensus.Range("D3").Value = 275
$wsCensus.Range("E3").Value = 0.06

# ---------------------------------------------------------------------
# Header / highlighted-number styling: 12pt black Calibri, applied to a
# seed cell via direct font formatting, then propagated to every other
# cell that needs it with a format-only copy/paste (keeps the style
# table minimal - matches a single extra font + a single extra cellXf).
# ---------------------------------------------------------------------
$wsCalf.Range("A1").Font.Size = 12
$wsCalf.Range("A1").Font.Color = 0

$wsCalf.Range("A1").Copy() | Out-Null
$wsCalf.Range("B1").PasteSpecial(-4122) | Out-Null
$wsCalf.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$wsCalf.Range("C2:E2").PasteSpecial(-4122) | Out-Null
$wsCensus.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$wsCensus.Range("C2:E2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row heights for the header + first data row on both sheets.
# ---------------------------------------------------------------------
$wsCalf.Rows.Item(1).RowHeight = 16
$wsCalf.Rows.Item(2).RowHeight = 16
$wsCensus.Rows.Item(1).RowHeight = 16
$wsCensus.Rows.Item(2).RowHeight = 16

# ---------------------------------------------------------------------
# Column widths on "proportion_calf" (best-fit-like widths carried over
# from the authored workbook).
# ---------------------------------------------------------------------
$wsCalf.Columns.Item(1).ColumnWidth = 17.5
$wsCalf.Columns.Item(2).ColumnWidth = 19.333333333333332
$wsCalf.Columns.Item(3).ColumnWidth = 16
$wsCalf.Columns.Item(4).ColumnWidth = 12.833333333333334
$wsCalf.Columns.Item(5).ColumnWidth = 15.666666666666666

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping, matching the authored workbook.
# ---------------------------------------------------------------------
$wsCalf.Range("B10").Select() | Out-Null
$wsCensus.Range("E13").Select() | Out-Null
$wsCensus.Activate() | Out-Null
